# Applies the "Adding search test cases" commit to the "Test Cases" sheet.
#
# Summary of changes:
#  1. Column D (Runmode) values for rows 2-106 change from "N" to "Y".
#  2. Row 106's D cell loses its stray fill/border style quirk (style 7 -> 3).
#  3. Row 107's B cell (TestCase_B106 / OPQA-1226) gains a border.
#  4. Two brand-new rows (108, 109) are appended with new test cases:
#       TestCase_B107 / OPQA-574 / "Verify that left navigation pane..."
#       TestCase_B108 / OPQA-569 / "Verify that sorting is retained..."
#  5. The active selection is moved to C122 to mirror the saved view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Flip Runmode column from N to Y for rows 2 through 106 -------------
$ws.Range("D2:D106").Value = "Y"

# --- 2. Normalize D106's formatting (drop the extra fill/border variant) ---
$ws.Range("D106").Interior.Pattern = -4142

# --- 3. Give B107 (OPQA-1226) a thin border like the rest of the table -----
$ws.Range("B107").Borders.LineStyle = 1

# --- 4. New row 108: TestCase_B107 ------------------------------------------
$ws.Range("A71").Copy()
$ws.Range("A108").PasteSpecial(-4122)
$ws.Range("B71").Copy()
$ws.Range("B108").PasteSpecial(-4122)
$ws.Range("C71").Copy()
$ws.Range("C108").PasteSpecial(-4122)
$ws.Range("D71").Copy()
$ws.Range("D108").PasteSpecial(-4122)
$ws.Range("E71").Copy()
$ws.Range("E108").PasteSpecial(-4122)

$ws.Range("A108").Value = "TestCase_B107"
$ws.Range("B108").Value = "OPQA-574"
$ws.Range("C108").Value = "Verify that left navigation pane content type is retained when user navigates back to ALL search results page from record view page"
$ws.Range("D108").Value = "Y"
$ws.Range("E108").Value = "SKIP"

# --- 5. New row 109: TestCase_B108 ------------------------------------------
$ws.Range("A71").Copy()
$ws.Range("A109").PasteSpecial(-4122)
$ws.Range("B71").Copy()
$ws.Range("B109").PasteSpecial(-4122)
$ws.Range("C71").Copy()
$ws.Range("C109").PasteSpecial(-4122)
$ws.Range("D71").Copy()
$ws.Range("D109").PasteSpecial(-4122)
$ws.Range("E71").Copy()
$ws.Range("E109").PasteSpecial(-4122)

$ws.Range("A109").Value = "TestCase_B108"
$ws.Range("B109").Value = "OPQA-569"
$ws.Range("C109").Value = "Verify that sorting is retained when user navigates back to ALL search results page from record view page"
$ws.Range("D109").Value = "Y"
$ws.Range("E109").Value = "PASS"

# --- 6. Restore the reported selection/active cell and scroll position ------
$ws.Activate()
$ws.Range("C122").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 106
$win.ScrollColumn = 1
